# Update the timestamp in column A (取得日時) for all data rows on the
# "ランサーズ" sheet from 2025-11-12 18:25:55 to 2025-11-12 18:33:37.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-11-12 18:25:55"
$newValue = "2025-11-12 18:33:37"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
